$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 147 holds the last existing record (date serial 45703) with a fixed
# set of values in columns B:J that continue unchanged into the new rows.
$lastRow = 147
$firstNewRow = 148
$lastNewRow = 153

# Copy formatting (date style/number format on column A) from the last
# existing row down into the new rows, then fill in the values.
$srcRange = $ws.Range("A$lastRow`:J$lastRow")
$destRange = $ws.Range("A$firstNewRow`:J$lastNewRow")
$srcRange.Copy()
$destRange.PasteSpecial(-4122) # xlPasteFormats

$startDate = $ws.Cells.Item($lastRow, 1).Value2

for ($r = $firstNewRow; $r -le $lastNewRow; $r++) {
    $offset = $r - $lastRow
    $ws.Cells.Item($r, 1).Value = $startDate + $offset
    $ws.Cells.Item($r, 2).Value = 116.4121952
    $ws.Cells.Item($r, 3).Value = 0.00170247
    $ws.Cells.Item($r, 4).Value = 0.008850780000000001
    $ws.Cells.Item($r, 5).Value = 0.06933635
    $ws.Cells.Item($r, 6).Value = 12792.90181321
    $ws.Cells.Item($r, 7).Value = 465.80531254
    $ws.Cells.Item($r, 8).Value = 0.24
    $ws.Cells.Item($r, 9).Value = 1.7904431
    $ws.Cells.Item($r, 10).Value = 485.38834923
}
